$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date that was bumped by one day
# (2023-09-15 -> 2023-09-16, i.e. serial 45184 -> 45185) for every data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value()
    if ($current -ne $null) {
        $cell.Value = $current.AddDays(1)
    }
}
